# [ADDITIONAL SCRAPING] add a "Player Info" sheet (ahead of the existing
# "ODI Batting" / "ODI Bowling" sheets) holding per-player identity fields,
# and replace the full match-scorecard URL column on the batting/bowling
# sheets with a bare numeric MATCH_CODE column.

$wb = $excel.ActiveWorkbook
$battingSheet = $wb.Sheets.Item("ODI Batting")

# Insert the new sheet before "ODI Batting" so the final tab order is
# Player Info, ODI Batting, ODI Bowling.
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# Writes $text into $cell as literal text, even when it looks numeric
# (e.g. "7128"), without leaving the cell's number format changed.
function Set-TextCell($cell, $text) {
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# Header row
$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $playerInfo.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

# Data row
Set-TextCell $playerInfo.Cells.Item(2, 1) "7128"
$playerInfo.Cells.Item(2, 2).Value = "Jack Jarvis"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Medium Fast"

# Rename MATCH_CARD_LINK -> MATCH_CODE and replace the full URL value with the
# bare match code on both the "ODI Batting" and "ODI Bowling" sheets.
foreach ($sheetName in @("ODI Batting", "ODI Bowling")) {
    $ws = $wb.Sheets.Item($sheetName)
    $found = $ws.Rows.Item(1).Find("MATCH_CARD_LINK")
    if ($found -ne $null) {
        $col = $found.Column
        $found.Value = "MATCH_CODE"
        Set-TextCell $ws.Cells.Item(2, $col) "4703"
    }
}
